$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Estadisticos 1P"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D2").Value = 15
$ws1.Range("F2").Value = 24
$ws1.Range("G2").Value = 61.54
$ws1.Range("H2").Value = 8

$ws1.Range("D3").Value = 16
$ws1.Range("F3").Value = 22
$ws1.Range("G3").Value = 57.89
$ws1.Range("H3").Value = 8

$ws1.Range("D4").Value = 10
$ws1.Range("F4").Value = 24
$ws1.Range("G4").Value = 70.59
$ws1.Range("H4").Value = 8.1

# ---------------------------------------------------------------
# Sheet "Estadisticos 2P"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value = 17
$ws2.Range("E2").Value = 2
$ws2.Range("F2").Value = 22
$ws2.Range("G2").Value = 56.41
$ws2.Range("H2").Value = 8.1

$ws2.Range("D3").Value = 20
$ws2.Range("E3").Value = 4
$ws2.Range("F3").Value = 18
$ws2.Range("G3").Value = 47.37
$ws2.Range("H3").Value = 8.3

$ws2.Range("D4").Value = 13
$ws2.Range("E4").Value = 3
$ws2.Range("F4").Value = 21
$ws2.Range("G4").Value = 61.76
$ws2.Range("H4").Value = 8.2

# ---------------------------------------------------------------
# Sheet "Estadisticos Final"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("D2").Value = 15
$ws3.Range("F2").Value = 24
$ws3.Range("G2").Value = 61.54

$ws3.Range("D3").Value = 16
$ws3.Range("F3").Value = 22
$ws3.Range("G3").Value = 57.89
$ws3.Range("H3").Value = 8.5

$ws3.Range("D4").Value = 10
$ws3.Range("F4").Value = 24
$ws3.Range("G4").Value = 70.59
$ws3.Range("H4").Value = 8.5

# ---------------------------------------------------------------
# Sheet "Rescatables" - add a new student row
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Range("A2").Value = 19330051920321
$ws4.Range("B2").Value = "DE JESUS"
$ws4.Range("C2").Value = "DE LA CRUZ"
$ws4.Range("D2").Value = "IGNACIO"
$ws4.Range("E2").Value = "ECOLOGÍA"
$ws4.Range("F2").Value = "4APV"
$ws4.Range("G2").Value = 2
